# Apply updated crypto price/volume data per the GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text type (matches original inlineStr cells) so numeric-looking
# strings like "214.00" or "1.00" keep their literal text, then restore
# the default "Normal" style so no stray formatting is introduced.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.982.06"
Set-TextValue $ws.Range("E2") "  +1.68%  "
Set-TextValue $ws.Range("D3") "1.648.05"
Set-TextValue $ws.Range("E3") "  +1.82%  "
Set-TextValue $ws.Range("E4") "  -0.09%  "
Set-TextValue $ws.Range("D5") "214.00"
Set-TextValue $ws.Range("E5") "  +1.42%  "
Set-TextValue $ws.Range("D6") "0.524"
Set-TextValue $ws.Range("E6") "  -0.32%  "
Set-TextValue $ws.Range("E7") "  -0.10%  "
Set-TextValue $ws.Range("D8") "23.71"
Set-TextValue $ws.Range("E8") "  +4.16%  "
Set-TextValue $ws.Range("E9") "  +1.93%  "
Set-TextValue $ws.Range("D10") "0.0614"
Set-TextValue $ws.Range("E10") "  +0.29%  "
Set-TextValue $ws.Range("E11") "  -1.67%  "
Set-TextValue $ws.Range("D12") "1.879.83"
Set-TextValue $ws.Range("D13") "1.650.56"
Set-TextValue $ws.Range("E13") "  +1.70%  "
Set-TextValue $ws.Range("E14") "  +1.50%  "
Set-TextValue $ws.Range("D15") "0.564"
Set-TextValue $ws.Range("E15") "  +2.54%  "
Set-TextValue $ws.Range("D16") "65.78"
Set-TextValue $ws.Range("E16") "  +1.85%  "
Set-TextValue $ws.Range("D17") "27.965.26"
Set-TextValue $ws.Range("E17") "  +1.60%  "
Set-TextValue $ws.Range("D18") "232.56"
Set-TextValue $ws.Range("E18") "  +1.58%  "
Set-TextValue $ws.Range("D19") "7.69"
Set-TextValue $ws.Range("E20") "  +0.70%  "
Set-TextValue $ws.Range("D21") "1.00"
Set-TextValue $ws.Range("E21") "  -0.08%  "
Set-TextValue $ws.Range("D22") "10.69"
Set-TextValue $ws.Range("E22") "  +7.29%  "
Set-TextValue $ws.Range("E23") "  +2.42%  "
Set-TextValue $ws.Range("D24") "2.17"
Set-TextValue $ws.Range("E24") "  +4.34%  "
Set-TextValue $ws.Range("D25") "151.46"
Set-TextValue $ws.Range("E25") "  +1.71%  "
Set-TextValue $ws.Range("D26") "6.94"
Set-TextValue $ws.Range("E26") "  +1.69%  "
Set-TextValue $ws.Range("D27") "15.74"
Set-TextValue $ws.Range("E27") "  +1.14%  "
Set-TextValue $ws.Range("E28") "  -0.03%  "
Set-TextValue $ws.Range("E29") "  -0.07%  "
Set-TextValue $ws.Range("E30") "  +1.52%  "
Set-TextValue $ws.Range("D31") "0.0485"
Set-TextValue $ws.Range("E32") "  +1.65%  "
Set-TextValue $ws.Range("D33") "1.455.76"
Set-TextValue $ws.Range("E33") "  +0.64%  "
Set-TextValue $ws.Range("E34") "  +2.06%  "
Set-TextValue $ws.Range("E35") "  +2.00%  "
Set-TextValue $ws.Range("E36") "  -0.63%  "
Set-TextValue $ws.Range("D37") "0.890"
Set-TextValue $ws.Range("E37") "  +3.28%  "
Set-TextValue $ws.Range("E38") "  +0.72%  "
Set-TextValue $ws.Range("E39") "  +0.74%  "
Set-TextValue $ws.Range("D40") "0.916"
Set-TextValue $ws.Range("E40") "  -2.94%  "
Set-TextValue $ws.Range("D41") "69.55"
Set-TextValue $ws.Range("E41") "  +0.05%  "
Set-TextValue $ws.Range("B42") "WEMIXToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D42") "1.02"
Set-TextValue $ws.Range("E42") "  +1.06%  "
Set-TextValue $ws.Range("B43") "PaxDollar"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D43") "1.00"
Set-TextValue $ws.Range("E43") "  +0.00%  "
Set-TextValue $ws.Range("E44") "  +0.34%  "
Set-TextValue $ws.Range("E45") "  +1.08%  "
Set-TextValue $ws.Range("E46") "  -0.22%  "
Set-TextValue $ws.Range("B47") "RenderToken"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D47") "1.78"
Set-TextValue $ws.Range("E47") "  +5.70%  "
Set-TextValue $ws.Range("B48") "RocketPoolETH"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D48") "1.788.45"
Set-TextValue $ws.Range("E48") "  +1.63%  "
Set-TextValue $ws.Range("D49") "88.79"
Set-TextValue $ws.Range("E49") "  +2.83%  "
Set-TextValue $ws.Range("B51") "EnergySwap"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "7.76"
Set-TextValue $ws.Range("E51") "  +0.16%  "
